$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4. This pushes the current row 4 data down to
# row 5, preserving all of its existing values/formatting untouched.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the updated observation record.
$ws.Range("A4").Value = 6760838
$ws.Range("B4").Value = 96239
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 504
$ws.Range("F4").Value = "Guckusko"
$ws.Range("G4").Value = "Cypripedium calceolus"
$ws.Range("H4").Value = "L."

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "146"

$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "överblommad"

$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "Södra Råda, vid Messmyrarnas SV-ände, Upl"

$ws.Range("Q4").Value = 713747.0846395431
$ws.Range("R4").Value = 6645995.92114509
$ws.Range("S4").Value = 50

$ws.Range("T4").Value = "Stockholm"
$ws.Range("U4").Value = "Norrtälje"
$ws.Range("V4").Value = "Uppland"
$ws.Range("W4").Value = "Söderby-Karl"

$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2006-06-26"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2006-06-26"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "00:00"

$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = "Endast blommande stjälkar räknade."

$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false

$ws.Range("AT4").NumberFormat = "@"
$ws.Range("AT4").Value = ""

$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value = "Gabriel Ekman"
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value = "Gabriel Ekman"
$ws.Range("AY4").NumberFormat = "@"
$ws.Range("AY4").Value = "Floraväkteri AB-län"
